# Replace the worksheet contents with the new, smaller data set.
# Original sheet held a B1:E1 header row plus A2:E3 sequence/activity data;
# the new sheet is just a one-column "sequence" list with a styled header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - remove all existing cell values and formatting.
$ws.Cells.Clear()

# Header cell (A1) plus two adjacent styled-but-empty header cells (B1:D1),
# matching the bold/bordered/centered "header" look used in the original file.
$header = $ws.Range("A1:D1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108  # xlCenter
$header.VerticalAlignment = -4160    # xlTop

$ws.Range("A1").Value = "sequence"

# Data rows - plain (unstyled) sequence values.
$ws.Range("A2").Value = "MVLT"
$ws.Range("A3").Value = "MVLQAHVELWNQT"
